$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 14-17 (Sending cluster = Resolving-Mac block)
$ws.Rows("14:17").Delete()

# Row 2: ECs -> ECs
$ws.Range("G2").Value = 0.9169710000000001
$ws.Range("H2").Value = 2.750913
$ws.Range("I2").Value = 0.01983441610607254
$ws.Range("J2").Value = 0.01983441610607253
$ws.Range("M2").Value = 0.4035516666666667
$ws.Range("N2").Value = 1.210655
$ws.Range("O2").Value = 0.1338129350090618
$ws.Range("P2").Value = 0.1338129350090617
$ws.Range("Q2").Value = 0.3700451753350001
$ws.Range("R2").Value = 3.330406578015
$ws.Range("S2").Value = 0.002654101433344572
$ws.Range("T2").Value = 0.002654101433344571

# Row 3: ECs -> FAPs
$ws.Range("G3").Value = 0.9169710000000001
$ws.Range("H3").Value = 2.750913
$ws.Range("I3").Value = 0.01983441610607254
$ws.Range("J3").Value = 0.01983441610607253
$ws.Range("M3").Value = 0.382402
$ws.Range("N3").Value = 1.147206
$ws.Range("O3").Value = 0.126799956981969
$ws.Range("P3").Value = 0.126799956981969
$ws.Range("Q3").Value = 0.350651544342
$ws.Range("R3").Value = 3.155863899078
$ws.Range("S3").Value = 0.002515003109012471
$ws.Range("T3").Value = 0.00251500310901247

# Row 4: ECs -> MuSCs
$ws.Range("G4").Value = 0.9169710000000001
$ws.Range("H4").Value = 2.750913
$ws.Range("I4").Value = 0.01983441610607254
$ws.Range("J4").Value = 0.01983441610607253
$ws.Range("M4").Value = 1.245432666666667
$ws.Range("N4").Value = 3.736298
$ws.Range("O4").Value = 0.4129706658366648
$ws.Range("P4").Value = 0.4129706658366648
$ws.Range("Q4").Value = 1.142025637786
$ws.Range("R4").Value = 10.278230740074
$ws.Range("S4").Value = 0.008191032025806243
$ws.Range("T4").Value = 0.008191032025806242

# Row 5: ECs -> Resolving-Mac
$ws.Range("G5").Value = 0.9169710000000001
$ws.Range("H5").Value = 2.750913
$ws.Range("I5").Value = 0.01983441610607254
$ws.Range("J5").Value = 0.01983441610607253
$ws.Range("M5").Value = 0.9844033333333333
$ws.Range("N5").Value = 2.95321
$ws.Range("O5").Value = 0.3264164421723045
$ws.Range("P5").Value = 0.3264164421723045
$ws.Range("Q5").Value = 0.9026693089700001
$ws.Range("R5").Value = 8.12402378073
$ws.Range("S5").Value = 0.00647427953790925
$ws.Range("T5").Value = 0.00647427953790925

# Row 6: FAPs -> ECs
$ws.Range("G6").Value = 40.163957
$ws.Range("H6").Value = 120.491871
$ws.Range("I6").Value = 0.8687609920100033
$ws.Range("J6").Value = 0.8687609920100032
$ws.Range("M6").Value = 0.4035516666666667
$ws.Range("N6").Value = 1.210655
$ws.Range("O6").Value = 0.1338129350090618
$ws.Range("P6").Value = 0.1338129350090617
$ws.Range("Q6").Value = 16.20823178727834
$ws.Range("R6").Value = 145.874086085505
$ws.Range("S6").Value = 0.1162514581622426
$ws.Range("T6").Value = 0.1162514581622426

# Row 7: FAPs -> FAPs
$ws.Range("G7").Value = 40.163957
$ws.Range("H7").Value = 120.491871
$ws.Range("I7").Value = 0.8687609920100033
$ws.Range("J7").Value = 0.8687609920100032
$ws.Range("M7").Value = 0.382402
$ws.Range("N7").Value = 1.147206
$ws.Range("O7").Value = 0.126799956981969
$ws.Range("P7").Value = 0.126799956981969
$ws.Range("Q7").Value = 15.358777484714
$ws.Range("R7").Value = 138.228997362426
$ws.Range("S7").Value = 0.1101588564144811
$ws.Range("T7").Value = 0.1101588564144811

# Row 8: FAPs -> MuSCs
$ws.Range("G8").Value = 40.163957
$ws.Range("H8").Value = 120.491871
$ws.Range("I8").Value = 0.8687609920100033
$ws.Range("J8").Value = 0.8687609920100032
$ws.Range("M8").Value = 1.245432666666667
$ws.Range("N8").Value = 3.736298
$ws.Range("O8").Value = 0.4129706658366648
$ws.Range("P8").Value = 0.4129706658366648
$ws.Range("Q8").Value = 50.02150407039534
$ws.Range("R8").Value = 450.193536633558
$ws.Range("S8").Value = 0.3587728053232925
$ws.Range("T8").Value = 0.3587728053232924

# Row 9: FAPs -> Resolving-Mac
$ws.Range("G9").Value = 40.163957
$ws.Range("H9").Value = 120.491871
$ws.Range("I9").Value = 0.8687609920100033
$ws.Range("J9").Value = 0.8687609920100032
$ws.Range("M9").Value = 0.9844033333333333
$ws.Range("N9").Value = 2.95321
$ws.Range("O9").Value = 0.3264164421723045
$ws.Range("P9").Value = 0.3264164421723045
$ws.Range("Q9").Value = 39.53753315065667
$ws.Range("R9").Value = 355.83779835591
$ws.Range("S9").Value = 0.2835778721099871
$ws.Range("T9").Value = 0.2835778721099871

# Row 10: MuSCs -> ECs
$ws.Range("G10").Value = 5.150379999999999
$ws.Range("H10").Value = 15.45114
$ws.Range("I10").Value = 0.1114045918839242
$ws.Range("J10").Value = 0.1114045918839242
$ws.Range("M10").Value = 0.4035516666666667
$ws.Range("N10").Value = 1.210655
$ws.Range("O10").Value = 0.1338129350090618
$ws.Range("P10").Value = 0.1338129350090617
$ws.Range("Q10").Value = 2.078444432966667
$ws.Range("R10").Value = 18.7059998967
$ws.Range("S10").Value = 0.0149073754134746
$ws.Range("T10").Value = 0.01490737541347459

# Row 11: MuSCs -> FAPs
$ws.Range("G11").Value = 5.150379999999999
$ws.Range("H11").Value = 15.45114
$ws.Range("I11").Value = 0.1114045918839242
$ws.Range("J11").Value = 0.1114045918839242
$ws.Range("M11").Value = 0.382402
$ws.Range("N11").Value = 1.147206
$ws.Range("O11").Value = 0.126799956981969
$ws.Range("P11").Value = 0.126799956981969
$ws.Range("Q11").Value = 1.96951561276
$ws.Range("R11").Value = 17.72564051484
$ws.Range("S11").Value = 0.0141260974584754
$ws.Range("T11").Value = 0.0141260974584754

# Row 12: MuSCs -> MuSCs
$ws.Range("G12").Value = 5.150379999999999
$ws.Range("H12").Value = 15.45114
$ws.Range("I12").Value = 0.1114045918839242
$ws.Range("J12").Value = 0.1114045918839242
$ws.Range("M12").Value = 1.245432666666667
$ws.Range("N12").Value = 3.736298
$ws.Range("O12").Value = 0.4129706658366648
$ws.Range("P12").Value = 0.4129706658366648
$ws.Range("Q12").Value = 6.414451497746666
$ws.Range("R12").Value = 57.73006347971999
$ws.Range("S12").Value = 0.04600682848756607
$ws.Range("T12").Value = 0.04600682848756607

# Row 13: MuSCs -> Resolving-Mac
$ws.Range("G13").Value = 5.150379999999999
$ws.Range("H13").Value = 15.45114
$ws.Range("I13").Value = 0.1114045918839242
$ws.Range("J13").Value = 0.1114045918839242
$ws.Range("M13").Value = 0.9844033333333333
$ws.Range("N13").Value = 2.95321
$ws.Range("O13").Value = 0.3264164421723045
$ws.Range("P13").Value = 0.3264164421723045
$ws.Range("Q13").Value = 5.070051239933332
$ws.Range("R13").Value = 45.63046115939999
$ws.Range("S13").Value = 0.03636429052440812
$ws.Range("T13").Value = 0.03636429052440812
